# Remove the first 9 FakeQuakes "sz" branch rows (b0-95 / b1-1 / b1-24 groups)
# from the "sz_weights_4_0_fq" worksheet, leaving only the 6 "hk" branch rows
# that follow (which shift up into rows 2:7).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sz_weights_4_0_fq")

$ws.Range("A2:A10").EntireRow.Delete()
